$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 251.375
$ws.Cells.Item(4, 9).Value = 295
$ws.Cells.Item(4, 10).Value = 225.2
$ws.Cells.Item(4, 11).Value = 295
$ws.Cells.Item(4, 12).Value = 225.2
$ws.Cells.Item(4, 13).Value = -181
$ws.Cells.Item(4, 14).Value = -453.2

$ws.Cells.Item(88, 8).Value = 1668.5294
$ws.Cells.Item(88, 9).Value = 1572.4166
$ws.Cells.Item(88, 10).Value = 1899.2
$ws.Cells.Item(88, 11).Value = 1572.4166
$ws.Cells.Item(88, 12).Value = 1899.2
$ws.Cells.Item(88, 13).Value = -1166.4166
$ws.Cells.Item(88, 14).Value = -2711.2

$ws.Cells.Item(91, 8).Value = 1668.5294
$ws.Cells.Item(91, 9).Value = 1572.4166
$ws.Cells.Item(91, 10).Value = 1899.2
$ws.Cells.Item(91, 11).Value = 1572.4166
$ws.Cells.Item(91, 12).Value = 1899.2
$ws.Cells.Item(91, 13).Value = -168.4166
$ws.Cells.Item(91, 14).Value = -4707.2

$ws.Cells.Item(94, 8).Value = 1911.5
$ws.Cells.Item(94, 9).Value = 1713.1428
$ws.Cells.Item(94, 11).Value = 1713.1428
$ws.Cells.Item(94, 13).Value = -1262.1428

$ws.Cells.Item(113, 8).Value = 6079.45
$ws.Cells.Item(113, 9).Value = 5823.375
$ws.Cells.Item(113, 11).Value = 5823.375
$ws.Cells.Item(113, 13).Value = -2569.375

$ws.Cells.Item(116, 8).Value = 2352
$ws.Cells.Item(116, 9).Value = 2084.25
$ws.Cells.Item(116, 10).Value = 2709
$ws.Cells.Item(116, 11).Value = 2084.25
$ws.Cells.Item(116, 12).Value = 2709
$ws.Cells.Item(116, 13).Value = 1357.75
$ws.Cells.Item(116, 14).Value = -9593

$ws.Cells.Item(121, 8).Value = 1133.5385
$ws.Cells.Item(121, 10).Value = 1133.5385
$ws.Cells.Item(121, 12).Value = 3400.6155
$ws.Cells.Item(121, 14).Value = -6894.6155

$ws.Cells.Item(132, 8).Value = 2161.0527
$ws.Cells.Item(132, 9).Value = 2004.875
$ws.Cells.Item(132, 10).Value = 2994
$ws.Cells.Item(132, 11).Value = 6014.625
$ws.Cells.Item(132, 12).Value = 8982
$ws.Cells.Item(132, 13).Value = -3484.625
$ws.Cells.Item(132, 14).Value = -14042

$ws.Cells.Item(137, 8).Value = 9205.107
$ws.Cells.Item(137, 9).Value = 1486.1428
$ws.Cells.Item(137, 10).Value = 16924.072
$ws.Cells.Item(137, 11).Value = 4458.428400000001
$ws.Cells.Item(137, 12).Value = 50772.216
$ws.Cells.Item(137, 13).Value = -1908.428400000001
$ws.Cells.Item(137, 14).Value = -55872.216

$ws.Cells.Item(138, 8).Value = 3985.75
$ws.Cells.Item(138, 10).Value = 4043.5322
$ws.Cells.Item(138, 12).Value = 12130.5966
$ws.Cells.Item(138, 14).Value = -22410.5966

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2464.7083
$ws.Cells.Item(61, 9).Value = 2007.85
$ws.Cells.Item(61, 11).Value = 2007.85
$ws.Cells.Item(61, 13).Value = -1795.85

$ws.Cells.Item(74, 8).Value = 2616.3914
$ws.Cells.Item(74, 9).Value = 1985
$ws.Cells.Item(74, 11).Value = 1985
$ws.Cells.Item(74, 13).Value = -1111

$ws.Cells.Item(77, 8).Value = 2616.3914
$ws.Cells.Item(77, 9).Value = 1985
$ws.Cells.Item(77, 11).Value = 9925
$ws.Cells.Item(77, 13).Value = -5557

$ws.Cells.Item(136, 8).Value = 2464.7083
$ws.Cells.Item(136, 9).Value = 2007.85
$ws.Cells.Item(136, 11).Value = 6023.549999999999
$ws.Cells.Item(136, 13).Value = -3473.549999999999

$ws.Cells.Item(139, 8).Value = 88333.336
$ws.Cells.Item(139, 10).Value = 88333.336
$ws.Cells.Item(139, 12).Value = 88333.336
$ws.Cells.Item(139, 14).Value = -98613.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 15799.8
$ws.Cells.Item(86, 9).Value = 15000
$ws.Cells.Item(86, 10).Value = 18999
$ws.Cells.Item(86, 11).Value = 15000
$ws.Cells.Item(86, 12).Value = 18999
$ws.Cells.Item(86, 13).Value = -13877
$ws.Cells.Item(86, 14).Value = -21245

$ws.Cells.Item(89, 8).Value = 15799.8
$ws.Cells.Item(89, 9).Value = 15000
$ws.Cells.Item(89, 10).Value = 18999
$ws.Cells.Item(89, 11).Value = 75000
$ws.Cells.Item(89, 12).Value = 94995
$ws.Cells.Item(89, 13).Value = -69384
$ws.Cells.Item(89, 14).Value = -106227

$ws.Cells.Item(105, 8).Value = 101368.5
$ws.Cells.Item(105, 9).Value = 1523.8125
$ws.Cells.Item(105, 10).Value = 500747.25
$ws.Cells.Item(105, 11).Value = 1523.8125
$ws.Cells.Item(105, 12).Value = 500747.25
$ws.Cells.Item(105, 13).Value = 223.1875
$ws.Cells.Item(105, 14).Value = -504241.25

$ws.Cells.Item(134, 8).Value = 4205.6445
$ws.Cells.Item(134, 9).Value = 4169.8604
$ws.Cells.Item(134, 10).Value = 4975
$ws.Cells.Item(134, 11).Value = 12509.5812
$ws.Cells.Item(134, 12).Value = 14925
$ws.Cells.Item(134, 13).Value = -9974.581199999999
$ws.Cells.Item(134, 14).Value = -19995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 906
$ws.Cells.Item(94, 9).Value = 865.6
$ws.Cells.Item(94, 11).Value = 865.6
$ws.Cells.Item(94, 13).Value = -414.6

$ws.Cells.Item(99, 8).Value = 2800
$ws.Cells.Item(99, 9).Value = 2210.6
$ws.Cells.Item(99, 11).Value = 2210.6
$ws.Cells.Item(99, 13).Value = -712.5999999999999

$ws.Cells.Item(105, 8).Value = 500869.25
$ws.Cells.Item(105, 9).Value = 1159
$ws.Cells.Item(105, 11).Value = 1159
$ws.Cells.Item(105, 13).Value = 588

$ws.Cells.Item(126, 8).Value = 2800
$ws.Cells.Item(126, 9).Value = 2210.6
$ws.Cells.Item(126, 11).Value = 6631.799999999999
$ws.Cells.Item(126, 13).Value = -4161.799999999999

$ws.Cells.Item(132, 8).Value = 2832.1365
$ws.Cells.Item(132, 9).Value = 2529
$ws.Cells.Item(132, 11).Value = 7587
$ws.Cells.Item(132, 13).Value = -5057

$ws.Cells.Item(138, 8).Value = 99989.39999999999
$ws.Cells.Item(138, 10).Value = 99989.39999999999
$ws.Cells.Item(138, 12).Value = 99989.39999999999
$ws.Cells.Item(138, 14).Value = -110269.4

$ws.Cells.Item(141, 8).Value = 271583.34
$ws.Cells.Item(141, 10).Value = 287272.72
$ws.Cells.Item(141, 12).Value = 287272.72
$ws.Cells.Item(141, 14).Value = -297632.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 96.27273
$ws.Cells.Item(10, 9).Value = 55.9
$ws.Cells.Item(10, 10).Value = 500
$ws.Cells.Item(10, 11).Value = 167.7
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = -28.69999999999999
$ws.Cells.Item(10, 14).Value = -1778

$ws.Cells.Item(29, 8).Value = 856.7143
$ws.Cells.Item(29, 9).Value = 499.25
$ws.Cells.Item(29, 10).Value = 1333.3334
$ws.Cells.Item(29, 11).Value = 1497.75
$ws.Cells.Item(29, 12).Value = 4000.0002
$ws.Cells.Item(29, 13).Value = -1220.75
$ws.Cells.Item(29, 14).Value = -4554.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 500000
$ws.Cells.Item(7, 9).Value = 500000
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 500000
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -499888
$ws.Cells.Item(7, 14).ClearContents()

$ws.Cells.Item(8, 8).Value = 500000
$ws.Cells.Item(8, 9).Value = 500000
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 500000
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -499861
$ws.Cells.Item(8, 14).ClearContents()

$ws.Cells.Item(14, 8).Value = 12535551
$ws.Cells.Item(14, 9).Value = 25000102
$ws.Cells.Item(14, 11).Value = 25000102
$ws.Cells.Item(14, 13).Value = -24999934

$ws.Cells.Item(34, 8).Value = 55499.5
$ws.Cells.Item(34, 10).Value = 55499.5
$ws.Cells.Item(34, 12).Value = 55499.5
$ws.Cells.Item(34, 14).Value = -56035.5

$ws.Cells.Item(39, 8).Value = 51991
$ws.Cells.Item(39, 10).Value = 51991
$ws.Cells.Item(39, 12).Value = 51991
$ws.Cells.Item(39, 14).Value = -53055

$ws.Cells.Item(76, 8).Value = 55499.5
$ws.Cells.Item(76, 10).Value = 55499.5
$ws.Cells.Item(76, 12).Value = 55499.5
$ws.Cells.Item(76, 14).Value = -56129.5

$ws.Cells.Item(79, 8).Value = 55499.5
$ws.Cells.Item(79, 10).Value = 55499.5
$ws.Cells.Item(79, 12).Value = 55499.5
$ws.Cells.Item(79, 14).Value = -57683.5

$ws.Cells.Item(80, 8).Value = 3111.2856
$ws.Cells.Item(80, 9).Value = 2666.3333
$ws.Cells.Item(80, 11).Value = 2666.3333
$ws.Cells.Item(80, 13).Value = -1668.3333

$ws.Cells.Item(83, 8).Value = 3111.2856
$ws.Cells.Item(83, 9).Value = 2666.3333
$ws.Cells.Item(83, 11).Value = 13331.6665
$ws.Cells.Item(83, 13).Value = -8339.666499999999

$ws.Cells.Item(102, 8).Value = 1895.6666
$ws.Cells.Item(102, 9).Value = 1895.6666
$ws.Cells.Item(102, 11).Value = 1895.6666
$ws.Cells.Item(102, 13).Value = -273.6666

$ws.Cells.Item(117, 8).Value = 33499.5
$ws.Cells.Item(117, 10).Value = 33499.5
$ws.Cells.Item(117, 12).Value = 33499.5
$ws.Cells.Item(117, 14).Value = -40383.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 62922.5
$ws.Cells.Item(22, 9).Value = 495
$ws.Cells.Item(22, 11).Value = 495
$ws.Cells.Item(22, 13).Value = -200

$ws.Cells.Item(27, 8).Value = 62922.5
$ws.Cells.Item(27, 9).Value = 495
$ws.Cells.Item(27, 11).Value = 495
$ws.Cells.Item(27, 13).Value = -388

$ws.Cells.Item(40, 8).Value = 5477.476
$ws.Cells.Item(40, 9).Value = 5151.6875
$ws.Cells.Item(40, 11).Value = 5151.6875
$ws.Cells.Item(40, 13).Value = -5015.6875

$ws.Cells.Item(93, 8).Value = 3024.606
$ws.Cells.Item(93, 9).Value = 2262.8262
$ws.Cells.Item(93, 11).Value = 2262.8262
$ws.Cells.Item(93, 13).Value = -1014.8262

$ws.Cells.Item(122, 8).Value = 42754.082
$ws.Cells.Item(122, 9).Value = 25104.9
$ws.Cells.Item(122, 11).Value = 75314.70000000001
$ws.Cells.Item(122, 13).Value = -72864.70000000001

$ws.Cells.Item(132, 8).Value = 4305
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 4305
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 12915
$ws.Cells.Item(132, 14).Value = -17975
$ws.Cells.Item(132, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2770.1667
$ws.Cells.Item(96, 9).Value = 2831.3333
$ws.Cells.Item(96, 11).Value = 2831.3333
$ws.Cells.Item(96, 13).Value = -1458.3333

$ws.Cells.Item(126, 8).Value = 2722
$ws.Cells.Item(126, 9).Value = 2919.125
$ws.Cells.Item(126, 11).Value = 8757.375
$ws.Cells.Item(126, 13).Value = -6287.375
